$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10178
$ws1.Range("F5").Value = 628

# Update "全部类型" (All Types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10178
$ws4.Range("F5").Value = 628
